# Update cryptocurrency price/volume figures per the Sat Feb 25 19:07:39 UTC 2023
# GitHub Actions refresh. Row 39/40 also swap coin identity (Algorand <-> InternetComputer).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.046.48'
$ws.Range("E2").Value = '  -0.40%  '
$ws.Range("D3").Value = '1.590.28'
$ws.Range("E3").Value = '  -0.63%  '
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("E5").Value = '  +0.05%  '
$ws.Range("D6").Value = "'301.70"
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("D7").Value = "'0.3768"
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").Value = "'0.3593"
$ws.Range("E8").Value = '  -1.85%  '
$ws.Range("D9").Value = "'50.77"
$ws.Range("E9").Value = '  +6.39%  '
$ws.Range("E10").Value = '  +0.06%  '
$ws.Range("E11").Value = '  -4.09%  '
$ws.Range("D12").Value = "'0.08061"
$ws.Range("E12").Value = '  -0.33%  '
$ws.Range("E13").Value = '  -4.25%  '
$ws.Range("D14").Value = "'6.487"
$ws.Range("D15").Value = "'7.263"
$ws.Range("E15").Value = '  -4.56%  '
$ws.Range("D16").Value = "'0.00001225"
$ws.Range("E16").Value = '  -3.43%  '
$ws.Range("D17").Value = '1.591.42'
$ws.Range("E17").Value = '  -0.26%  '
$ws.Range("D18").Value = "'92.70"
$ws.Range("E18").Value = '  +1.16%  '
$ws.Range("D19").Value = "'0.06829"
$ws.Range("E19").Value = '  +0.44%  '
$ws.Range("D20").Value = "'17.95"
$ws.Range("E20").Value = '  -2.70%  '
$ws.Range("D21").Value = "'6.458"
$ws.Range("E21").Value = '  -2.09%  '
$ws.Range("D22").Value = "'1.003"
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("D23").Value = "'12.87"
$ws.Range("E23").Value = '  -1.51%  '
$ws.Range("D24").Value = '23.051.38'
$ws.Range("E24").Value = '  -0.43%  '
$ws.Range("D25").Value = "'2.372"
$ws.Range("E25").Value = '  +0.28%  '
$ws.Range("D26").Value = "'2.817"
$ws.Range("E26").Value = '  -3.01%  '
$ws.Range("D27").Value = "'20.88"
$ws.Range("E27").Value = '  -1.15%  '
$ws.Range("D28").Value = "'148.49"
$ws.Range("E28").Value = '  -1.83%  '
$ws.Range("D29").Value = "'5.224"
$ws.Range("E29").Value = '  -0.35%  '
$ws.Range("D30").Value = "'133.44"
$ws.Range("E30").Value = '  +1.19%  '
$ws.Range("D31").Value = "'2.358"
$ws.Range("E31").Value = '  -3.75%  '
$ws.Range("D32").Value = "'6.518"
$ws.Range("E32").Value = '  -8.19%  '
$ws.Range("D33").Value = '1.767.25'
$ws.Range("E33").Value = '  -0.27%  '
$ws.Range("D34").Value = "'0.9429"
$ws.Range("E34").Value = '  -3.86%  '
$ws.Range("D35").Value = "'0.07399"
$ws.Range("E35").Value = '  -4.48%  '
$ws.Range("D36").Value = "'10.12"
$ws.Range("E36").Value = '  +0.38%  '
$ws.Range("D37").Value = "'0.02678"
$ws.Range("E37").Value = '  -3.79%  '
$ws.Range("D38").Value = "'0.08772"
$ws.Range("E38").Value = '  -1.18%  '
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").Value = "'0.2478"
$ws.Range("E39").Value = '  -2.62%  '
$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").Value = "'6.044"
$ws.Range("E40").Value = '  -4.31%  '
$ws.Range("D41").Value = "'1.347"
$ws.Range("E41").Value = '  -3.56%  '
$ws.Range("D42").Value = "'0.6921"
$ws.Range("E42").Value = '  -3.53%  '
$ws.Range("D43").Value = "'12.15"
$ws.Range("E43").Value = '  -4.97%  '
$ws.Range("D44").Value = "'14.85"
$ws.Range("E44").Value = '  -6.60%  '
$ws.Range("D45").Value = "'0.6462"
$ws.Range("E45").Value = '  -2.89%  '
$ws.Range("D46").Value = "'4.002"
$ws.Range("E46").Value = '  +0.80%  '
$ws.Range("D47").Value = "'2.252"
$ws.Range("E47").Value = '  -2.72%  '
$ws.Range("D48").Value = "'131.36"
$ws.Range("E48").Value = '  -0.69%  '
$ws.Range("D49").Value = "'0.07895"
$ws.Range("E49").Value = '  -1.01%  '
$ws.Range("D50").Value = "'1.197"
$ws.Range("E50").Value = '  +1.94%  '
$ws.Range("D51").Value = "'1.213"
$ws.Range("E51").Value = '  +3.98%  '
